{"js": "// Change: the first paragraph (\"El objetivo es aprender lo que hace git\")\n// becomes \"jejejej\", and a brand-new paragraph carrying the original text\n// (split into two runs around \"git\", wrapped with proofErr spell-check\n// markers, exactly like Word emits when it re-checks spelling) is inserted\n// right after it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstPara = paragraphs.items[0];\n\n// Insert the new paragraph (with the original wording) right after the\n// first paragraph, using raw OOXML so we can reproduce the exact\n// proofErr-wrapped run split that appears in the target markup.\nconst endOfFirstPara = firstPara.getRange(\"End\");\n\nconst newParagraphOoxml = `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:ind w:left=\"1416\"/>\n              <w:rPr>\n                <w:rFonts w:asciiTheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"28\"/>\n                <w:szCs w:val=\"28\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:asciiTheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"28\"/>\n                <w:szCs w:val=\"28\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">El objetivo es aprender lo que hace </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:asciiTheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"28\"/>\n                <w:szCs w:val=\"28\"/>\n              </w:rPr>\n              <w:t>git</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nendOfFirstPara.insertOoxml(newParagraphOoxml, Word.InsertLocation.after);\n\n// Replace the first paragraph's own text with \"jejejej\".\nfirstPara.insertText(\"jejejej\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Change: the first paragraph (\"El objetivo es aprender lo que hace git\")\n# becomes \"jejejej\", and a brand-new paragraph carrying the original\n# wording (split into two runs around \"git\", wrapped with proofErr\n# spell-check markers, exactly like Word emits when it re-checks\n# spelling) is inserted right after it.\n\n$d = $word.ActiveDocument\n$firstPara = $d.Paragraphs.Item(1)\n\n# Create an empty paragraph right after the first one; it inherits the\n# same paragraph/run formatting (indent 1416, minorHAnsi theme fonts,\n# size 28) because it is split off of the first paragraph's mark.\n$firstPara.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs.Item(2)\n\n$newParagraphOoxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:ind w:left=\"1416\"/>\n              <w:rPr>\n                <w:rFonts w:asciiTheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"28\"/>\n                <w:szCs w:val=\"28\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:asciiTheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"28\"/>\n                <w:szCs w:val=\"28\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">El objetivo es aprender lo que hace </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:asciiTheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"28\"/>\n                <w:szCs w:val=\"28\"/>\n              </w:rPr>\n              <w:t>git</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$newPara.Range.InsertXML($newParagraphOoxml)\n\n# Replace the first paragraph's own text with \"jejejej\".\n$firstPara.Range.Text = \"jejejej\"\n"}
